# edit.ps1 - applies the PLANTILLA_UNT "objetivos" block expansion plus the
# related lastRenderedPageBreak move (Material/Metodos gains the break that
# Operacionalizacion used to carry).
$d = $word.ActiveDocument

function Find-ParagraphIndexByExactText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

function Find-ParagraphIndexByContains($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

# 1) ${objetivos} -> the block_txt_obj_general / block_obj_general / ... tree
$objIdx = Find-ParagraphIndexByExactText $d '${objetivos}'
if ($objIdx -lt 0) {
    throw 'Could not locate the ${objetivos} paragraph'
}
$objXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block_txt_obj_general</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t>General:</w:t></w:r></w:p><w:p><w:r><w:t>${/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block_txt_obj_general</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block_obj_general</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>obj_descripcion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>block_obj_general</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>block_txt_obj_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>especifico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Espec</w:t></w:r><w:r><w:t>í</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>fico</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>${/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block_txt_obj_</w:t></w:r><w:r><w:t>especifico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>block_obj_especifico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>obj_descripcion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>block_obj_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>especifico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>'
$d.Paragraphs.Item($objIdx).Range.InsertXML($objXml)
Write-Host "Replaced objetivos placeholder at paragraph $objIdx"

# 2) "Material, Metodos y Tecnicas" gains a <w:lastRenderedPageBreak/>
$matIdx = Find-ParagraphIndexByContains $d "todos y T"
if ($matIdx -lt 0) {
    throw 'Could not locate the Material, Metodos y Tecnicas paragraph'
}
$matXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo3"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Material, Métodos y Técnicas</w:t></w:r></w:p>'
$d.Paragraphs.Item($matIdx).Range.InsertXML($matXml)
Write-Host "Added lastRenderedPageBreak to paragraph $matIdx"

# 3) "Operacionalizacion de Variables..." loses its <w:lastRenderedPageBreak/>
$opIdx = Find-ParagraphIndexByContains $d "matriz de consistencia"
if ($opIdx -lt 0) {
    throw 'Could not locate the Operacionalizacion paragraph'
}
$opXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo3"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Operacionalización de Variables y matriz de consistencia</w:t></w:r></w:p>'
$d.Paragraphs.Item($opIdx).Range.InsertXML($opXml)
Write-Host "Removed lastRenderedPageBreak from paragraph $opIdx"

Write-Host "Done"
